# Implementa suporte completo para PT-BR no template de Prompts:
# traduz valores de enum (categoria, tipo, idioma, comportamento, tom,
# dificuldade) para português e ajusta larguras de algumas colunas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Larguras de coluna (B, C, E ficam mais largas) ---
$ws.Columns.Item(2).ColumnWidth = 20
$ws.Columns.Item(3).ColumnWidth = 20
$ws.Columns.Item(5).ColumnWidth = 15

# --- Linha 2: Pergunta sobre experiência ---
$ws.Range("B2").Value = "Entrevista"
$ws.Range("C2").Value = "Pergunta de Acompanhamento"
$ws.Range("E2").Value = "Português (Brasil)"
$ws.Range("F2").Value = "Profissional"
$ws.Range("G2").Value = "Neutro"
$ws.Range("K2").Value = "Médio"
$ws.Range("M2").Value = "experiência, atendimento, cliente"

# --- Linha 3: Avaliação de soft skills ---
$ws.Range("B3").Value = "Entrevista"
$ws.Range("C3").Value = "Avaliação"
$ws.Range("E3").Value = "Português (Brasil)"
$ws.Range("F3").Value = "Profissional"
$ws.Range("G3").Value = "Encorajador"
$ws.Range("K3").Value = "Difícil"

# --- Linha 4: Atendimento reativo ---
$ws.Range("B4").Value = "Atendimento ao Cliente"
$ws.Range("C4").Value = "Mensagem Inicial"
$ws.Range("E4").Value = "Português (Brasil)"
$ws.Range("F4").Value = "Amigável"
$ws.Range("G4").Value = "Positivo"
$ws.Range("K4").Value = "Médio"

# --- Linha 5: Venda consultiva ---
$ws.Range("B5").Value = "Vendas"
$ws.Range("C5").Value = "Pergunta de Acompanhamento"
$ws.Range("E5").Value = "Português (Brasil)"
$ws.Range("F5").Value = "Profissional"
$ws.Range("G5").Value = "Desafiador"
$ws.Range("K5").Value = "Especialista"
